$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Summary  (B3, B4, B5, B6, B7, B9)
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1399.77    # Current Capital
$summary.Range("B4").Value = -0.24      # Total P&L $
$summary.Range("B5").Value = -0.3       # Total P&L %
$summary.Range("B6").Value = 16         # Total Trades
$summary.Range("B7").Value = 7          # Winning Trades
$summary.Range("B9").Value = 43.75      # Win Rate %

# ---------------------------------------------------------------------------
# Sheet: Strategy Status  (C5:G5 - MarketMaking row)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 99.77   # Capital
$status.Range("D5").Value = 16      # Trades
$status.Range("E5").Value = -0.24   # P&L $
$status.Range("F5").Value = -0.23   # P&L %
$status.Range("G5").Value = 43.75   # Win Rate %

# ---------------------------------------------------------------------------
# Sheets: "All Trades" and "MarketMaking" - append new Trade #16 row (row 17)
# ---------------------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Cells.Item(17, 1).Value = 16                 # A17 Trade #

    # Keep Date/Time as plain text, not auto-converted to date/time serials
    $ws.Cells.Item(17, 2).NumberFormat = "@"
    $ws.Cells.Item(17, 2).Value = "2026-02-17"        # B17 Date
    $ws.Cells.Item(17, 3).NumberFormat = "@"
    $ws.Cells.Item(17, 3).Value = "20:03:42"          # C17 Time

    $ws.Cells.Item(17, 4).Value = "MarketMaking"      # D17 Strategy
    $ws.Cells.Item(17, 5).Value = "UP"                # E17 Side
    $ws.Cells.Item(17, 6).Value = 0.87                # F17 Entry Price
    $ws.Cells.Item(17, 7).Value = 0.93                # G17 Exit Price
    $ws.Cells.Item(17, 8).Value = "CLOSED"            # H17 Status
    $ws.Cells.Item(17, 9).Value = 6.8966              # I17 P&L %
    $ws.Cells.Item(17, 10).Value = 0.06               # J17 P&L $
    $ws.Cells.Item(17, 11).Value = 99.77              # K17 Capital After
    $ws.Cells.Item(17, 12).Value = 0                  # L17 Entry Slippage (bps)
    $ws.Cells.Item(17, 13).Value = 0                  # M17 Exit Slippage (bps)
    $ws.Cells.Item(17, 14).Value = 0.6                # N17 Confidence
    $ws.Cells.Item(17, 15).Value = "Normal spread capture: 19600 bps"  # O17 Entry Reason
    $ws.Cells.Item(17, 16).Value = "early_exit"       # P17 Exit Reason
    $ws.Cells.Item(17, 17).Value = 0.13               # Q17 Duration (min)
}
